$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of A3 (removes the cell entirely) and B3 (keep formatting/style, clear value)
$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()

# Update selection to match the diff: active cell A3, selected range A3:B3
$ws.Range("A3:B3").Select()

$wb.Save()
